$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9855021238327026
$ws.Range("B1").Value = 1.980955481529236
$ws.Range("C1").Value = 8.750088691711426
$ws.Range("D1").Value = 2.800483703613281
$ws.Range("E1").Value = 1.425874590873718
